$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.804.65'
$ws.Range('E2').Value = '  +1.33%  '

$ws.Range('D3').Value = '3.910.59'
$ws.Range('E3').Value = '  +1.10%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '''607.77'
$ws.Range('E5').Value = '  +1.10%  '

$ws.Range('D6').Value = '''169.74'
$ws.Range('E6').Value = '  +4.35%  '

$ws.Range('D7').Value = '3.911.08'
$ws.Range('E7').Value = '  +1.17%  '

$ws.Range('E8').Value = '  +0.23%  '

$ws.Range('D9').Value = '''0.536'
$ws.Range('E9').Value = '  +0.89%  '

$ws.Range('E10').Value = '  +0.61%  '

$ws.Range('E11').Value = '  +1.51%  '

$ws.Range('D12').Value = '''0.469'
$ws.Range('E12').Value = '  +2.20%  '

$ws.Range('B13').Value = 'ShibaInu'
$ws.Range('C13').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D13').Value = '''0.0000256'
$ws.Range('E13').Value = '  +4.87%  '

$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = '''38.29'
$ws.Range('E14').Value = '  +3.44%  '

$ws.Range('D15').Value = '4.573.86'
$ws.Range('E15').Value = '  +1.32%  '

$ws.Range('D16').Value = '3.905.17'
$ws.Range('E16').Value = '  +0.99%  '

$ws.Range('D17').Value = '69.852.83'
$ws.Range('E17').Value = '  +1.13%  '

$ws.Range('D18').Value = '''18.79'
$ws.Range('E18').Value = '  +9.96%  '

$ws.Range('D19').Value = '''7.63'
$ws.Range('E19').Value = '  +0.90%  '

$ws.Range('E20').Value = '  -0.85%  '

$ws.Range('D21').Value = '''11.21'
$ws.Range('E21').Value = '  -1.47%  '

$ws.Range('D22').Value = '''493.23'
$ws.Range('E22').Value = '  +1.66%  '

$ws.Range('E23').Value = '  +3.89%  '

$ws.Range('E24').Value = '  +3.50%  '

$ws.Range('D25').Value = '''85.63'
$ws.Range('E25').Value = '  +1.93%  '

$ws.Range('D26').Value = '''2.30'
$ws.Range('E26').Value = '  +2.45%  '

$ws.Range('E27').Value = '  +2.23%  '

$ws.Range('D28').Value = '''10.18'
$ws.Range('E28').Value = '  +2.04%  '

$ws.Range('E29').Value = '  +0.10%  '

$ws.Range('E30').Value = '  +1.25%  '

$ws.Range('D31').Value = '4.064.77'
$ws.Range('E31').Value = '  +1.22%  '

$ws.Range('E32').Value = '  +2.75%  '

$ws.Range('D33').Value = '''7.84'
$ws.Range('E33').Value = '  -1.17%  '

$ws.Range('D34').Value = '''32.11'
$ws.Range('E34').Value = '  -0.80%  '

$ws.Range('D35').Value = '3.877.62'
$ws.Range('E35').Value = '  +1.66%  '

$ws.Range('D36').Value = '''0.107'
$ws.Range('E36').Value = '  +0.59%  '

$ws.Range('D37').Value = '''6.13'
$ws.Range('E37').Value = '  +3.89%  '

$ws.Range('D38').Value = '''1.04'
$ws.Range('E38').Value = '  +1.23%  '

$ws.Range('E39').Value = '  +1.47%  '

$ws.Range('D40').Value = '''3.31'
$ws.Range('E40').Value = '  +11.54%  '

$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.11%  '

$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').Value = '''0.330'
$ws.Range('E42').Value = '  +3.38%  '

$ws.Range('D43').Value = '''2.13'
$ws.Range('E43').Value = '  +7.12%  '

$ws.Range('D44').Value = '''438.80'
$ws.Range('E44').Value = '  +0.35%  '

$ws.Range('D45').Value = '''48.17'
$ws.Range('E45').Value = '  -0.68%  '

$ws.Range('D46').Value = '''8.69'
$ws.Range('E46').Value = '  +3.36%  '

$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('E48').Value = '  +3.12%  '

$ws.Range('D49').Value = '''40.59'
$ws.Range('E49').Value = '  +4.59%  '

$ws.Range('D50').Value = '''143.71'
$ws.Range('E50').Value = '  +0.17%  '

$ws.Range('D51').Value = '''0.000271'
$ws.Range('E51').Value = '  +18.88%  '
